$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value = "2019-2020"
$ws.Range("P4").Value = "FIDELITY"
$ws.Range("M20").Value = "2"
$ws.Range("O20").Value = "FUMINO ONA FURAHASHI "
$ws.Range("E7").Value = "13"
$ws.Range("S7").Value = "Yolo 3_7"

$ws.Range("A8").Value = "123543457474"
$ws.Range("B8").Value = "Paderogao, Phil Rey, E. Jr"
$ws.Range("C8").Value = "F"
$ws.Range("D8").Value = "2006-12-11"
$ws.Range("E8").Value = "12"
$ws.Range("F8").Value = "Manay, Davao Oriental"
$ws.Range("G8").Value = "Mandaya"
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = "Roman Catholic"
$ws.Range("J8").Value = "Purok 18, Bato St."
$ws.Range("K8").Value = "Central"
$ws.Range("L8").Value = "Manay"
$ws.Range("M8").Value = "Davao Oriental"
$ws.Range("N8").Value = "Enrique C. Paderogao"
$ws.Range("O8").Value = "Jocelyn E. Paderogao"
$ws.Range("P8").Value = ""
$ws.Range("Q8").Value = ""
$ws.Range("R8").Value = "09483428056"
$ws.Range("S8").Value = "YoloNew 1_7"
